# Update "想去人数" (want-to-go count) figures in column F on the two
# sheets that carry the full event listing ("展览" and "全部类型").
# "演出" / "本地生活" only contain a header row and are untouched.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F3"  = 3180
    "F8"  = 1630
    "F13" = 191
    "F15" = 230
    "F16" = 239
    "F17" = 231
    "F23" = 377
    "F24" = 206
    "F25" = 103
    "F28" = 24
    "F29" = 229
    "F30" = 2156
    "F31" = 8
    "F34" = 323
    "F38" = 344
    "F40" = 514
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
